$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text to be R-friendly (lowercase, underscore instead of hyphen)
$ws.Range("A1").Value = "usual"
$ws.Range("B1").Value = "slow_release"

# Move the selection/active cell to B1 (matches the saved cursor position)
$ws.Range("B1").Select()
